# Updates the crypto price/volume table (columns D and E) to the latest
# scraped values, as produced by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D, omit if unchanged), new Volume(1h) (column E, omit if unchanged)
$updates = @(
    @{ Row = 2;  D = '28.025.01';  E = '  -0.23%  ' },
    @{ Row = 3;  D = '1.869.17';   E = '  -1.21%  ' },
    @{ Row = 4;                    E = '  +0.43%  ' },
    @{ Row = 5;  D = '312.58';     E = '  -0.49%  ' },
    @{ Row = 6;                    E = '  +0.33%  ' },
    @{ Row = 7;  D = '0.5059';     E = '  +0.32%  ' },
    @{ Row = 8;  D = '0.3801';     E = '  -2.55%  ' },
    @{ Row = 9;  D = '0.08296';    E = '  -10.26%  ' },
    @{ Row = 10; D = '1.110';      E = '  -1.63%  ' },
    @{ Row = 11; D = '41.51';      E = '  -0.72%  ' },
    @{ Row = 12; D = '6.197';      E = '  -2.91%  ' },
    @{ Row = 13; D = '1.877.29';   E = '  -0.76%  ' },
    @{ Row = 14; D = '20.46';      E = '  -1.73%  ' },
    @{ Row = 15; D = '7.188';      E = '  -1.47%  ' },
    @{ Row = 16; D = '1.004';      E = '  +0.35%  ' },
    @{ Row = 17; D = '0.00001094'; E = '  -1.19%  ' },
    @{ Row = 18; D = '90.72';      E = '  -1.58%  ' },
    @{ Row = 19; D = '0.06630';    E = '  -0.29%  ' },
    @{ Row = 20; D = '17.92';      E = '  +0.44%  ' },
    @{ Row = 22; D = '6.031';      E = '  -2.82%  ' },
    @{ Row = 23; D = '28.073.39';  E = '  -0.23%  ' },
    @{ Row = 24; D = '11.14';      E = '  -2.27%  ' },
    @{ Row = 25; D = '2.265';      E = '  -2.31%  ' },
    @{ Row = 26; D = '2.584';      E = '  +1.68%  ' },
    @{ Row = 27; D = '2.088.52';   E = '  -0.98%  ' },
    @{ Row = 28; D = '156.93';     E = '  -0.88%  ' },
    @{ Row = 29; D = '20.51';      E = '  -1.45%  ' },
    @{ Row = 30; D = '125.65';     E = '  -0.96%  ' },
    @{ Row = 31; D = '0.1053';     E = '  -0.15%  ' },
    @{ Row = 32; D = '1.043';      E = '  -3.34%  ' },
    @{ Row = 33; D = '5.595';      E = '  -0.16%  ' },
    @{ Row = 34; D = '3.604';      E = '  +0.09%  ' },
    @{ Row = 35; D = '9.726';      E = '  +2.75%  ' },
    @{ Row = 36; D = '0.02444';    E = '  +1.68%  ' },
    @{ Row = 37; D = '0.06579';    E = '  -0.55%  ' },
    @{ Row = 38; D = '0.2165';     E = '  -1.45%  ' },
    @{ Row = 39; D = '1.210';      E = '  -0.67%  ' },
    @{ Row = 40; D = '0.6469';     E = '  +0.35%  ' },
    @{ Row = 41; D = '1.243';      E = '  -7.37%  ' },
    @{ Row = 42; D = '11.29';      E = '  -2.78%  ' },
    @{ Row = 43; D = '4.876';      E = '  -1.90%  ' },
    @{ Row = 44; D = '0.6123';     E = '  +1.13%  ' },
    @{ Row = 45; D = '13.10';      E = '  -1.41%  ' },
    @{ Row = 46; D = '1.295';      E = '  -0.53%  ' },
    @{ Row = 47; D = '3.664';      E = '  -0.64%  ' },
    @{ Row = 48; D = '2.008';      E = '  +0.25%  ' },
    @{ Row = 49; D = '1.213';      E = '  +1.57%  ' },
    @{ Row = 50; D = '121.18';     E = '  -0.66%  ' },
    @{ Row = 51; D = '80.17';      E = '  +1.65%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey('D')) {
        # Prefix with an apostrophe so Excel stores the value as literal text
        # instead of re-interpreting strings like "28.025.01" or "0.08296" as
        # a number/date and silently reformatting / truncating them.
        $ws.Range("D" + $u.Row).Value = "'" + $u.D
        # Re-apply the plain "Normal" style so the quote-prefix formatting
        # introduced above doesn't leave a stray cell style behind.
        $ws.Range("D" + $u.Row).Style = "Normal"
    }
    if ($u.ContainsKey('E')) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
